# Refresh cryptos list: updates Price (col D) and Volume(1h) (col E) cells
# for rows 2-51. Numeric-looking prices (single-dot decimals) are written
# with a leading apostrophe so Excel keeps them as text (matching the
# original inlineStr storage) instead of silently re-typing them as
# numbers and dropping significant trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.998.33'
$ws.Range("E2").Value = '  -0.12%  '
$ws.Range("D3").Value = '1.677.07'
$ws.Range("E3").Value = '  +0.24%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '''215.17'
$ws.Range("E5").Value = '  -0.37%  '
$ws.Range("D6").Value = '''0.520'
$ws.Range("E6").Value = '  +1.63%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("E8").Value = '  +0.13%  '
$ws.Range("E9").Value = '  +0.30%  '
$ws.Range("D10").Value = '''20.33'
$ws.Range("E10").Value = '  +0.98%  '
$ws.Range("D11").Value = '''0.0887'
$ws.Range("E11").Value = '  -0.39%  '
$ws.Range("D12").Value = '1.912.79'
$ws.Range("E12").Value = '  +0.15%  '
$ws.Range("D13").Value = '1.669.05'
$ws.Range("E13").Value = '  -0.27%  '
$ws.Range("E14").Value = '  +0.20%  '
$ws.Range("D15").Value = '''0.529'
$ws.Range("E15").Value = '  +1.47%  '
$ws.Range("D16").Value = '''65.78'
$ws.Range("E16").Value = '  -0.04%  '
$ws.Range("D17").Value = '27.021.57'
$ws.Range("E17").Value = '  -0.13%  '
$ws.Range("D18").Value = '''8.15'
$ws.Range("E18").Value = '  +5.78%  '
$ws.Range("D19").Value = '''236.88'
$ws.Range("E19").Value = '  +0.82%  '
$ws.Range("D20").Value = '0.0₃0734'
$ws.Range("E20").Value = '  -0.42%  '
$ws.Range("E21").Value = '  +0.09%  '
$ws.Range("E22").Value = '  -0.33%  '
$ws.Range("E23").Value = '  -0.66%  '
$ws.Range("E24").Value = '  -2.45%  '
$ws.Range("D25").Value = '''146.11'
$ws.Range("E25").Value = '  +0.56%  '
$ws.Range("D26").Value = '''7.23'
$ws.Range("E26").Value = '  +0.83%  '
$ws.Range("D27").Value = '''16.14'
$ws.Range("E27").Value = '  +1.48%  '
$ws.Range("E28").Value = '  -1.51%  '
$ws.Range("E29").Value = '  +0.01%  '
$ws.Range("D30").Value = '''0.0498'
$ws.Range("E30").Value = '  +0.08%  '
$ws.Range("E31").Value = '  -0.38%  '
$ws.Range("E32").Value = '  -0.29%  '
$ws.Range("D33").Value = '1.479.15'
$ws.Range("E33").Value = '  +1.79%  '
$ws.Range("E34").Value = '  +0.62%  '
$ws.Range("E35").Value = '  +5.08%  '
$ws.Range("E36").Value = '  +0.32%  '
$ws.Range("D37").Value = '''0.582'
$ws.Range("E37").Value = '  +2.15%  '
$ws.Range("E38").Value = '  +2.57%  '
$ws.Range("E39").Value = '  +1.28%  '
$ws.Range("D40").Value = '''5.87'
$ws.Range("E40").Value = '  -3.50%  '
$ws.Range("E41").Value = '  +1.56%  '
$ws.Range("E42").Value = '  +0.09%  '
$ws.Range("E43").Value = '  +1.95%  '
$ws.Range("D44").Value = '''67.40'
$ws.Range("E44").Value = '  +2.36%  '
$ws.Range("D45").Value = '1.822.01'
$ws.Range("E45").Value = '  +0.08%  '
$ws.Range("E46").Value = '  +0.16%  '
$ws.Range("D47").Value = '''90.45'
$ws.Range("E47").Value = '  +0.16%  '
$ws.Range("E48").Value = '  +0.79%  '
$ws.Range("E49").Value = '  -0.37%  '
$ws.Range("E50").Value = '  +1.33%  '
$ws.Range("E51").Value = '  -0.22%  '
